# "Idea for weighting detection rates"
#
# Inserts a new constant (program_prop_smearnegextrapul_diagnosis, =2/3)
# above the existing "program_prop_death_reporting" row on the
# "constants" sheet, and updates the active sheet/selection bookkeeping
# to match (constants becomes the active tab, with B29 selected; the
# previously-active time_variants tab is no longer marked selected).

$wb = $excel.ActiveWorkbook

$wsConst = $wb.Worksheets.Item("constants")

# Make "constants" the active sheet (this also clears tabSelected on the
# sheet that was previously active, i.e. "time_variants").
$wsConst.Activate()

# Insert a fresh row above the current row 42 ("program_prop_death_reporting"),
# pushing it (and everything below) down by one. Excel carries the
# formatting of the row above into the new row automatically.
$wsConst.Rows.Item(42).Insert() | Out-Null

# Populate the new row with the new constant.
$wsConst.Range("A42").Value = "program_prop_smearnegextrapul_diagnosis"
$wsConst.Range("B42").Formula = "=2/3"

# Match the final selection recorded for the "constants" sheet.
$null = $wsConst.Range("B29").Select()
